$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 (Mod/Dense SAV under Fall season) - fill in new H:L letter cells
$ws.Range("H7").Value = "a"
$ws.Range("I7").Value = "b"
$ws.Range("J7").Value = "a"
$ws.Range("K7").Value = "b"
$ws.Range("L7").Value = "a"

# Row 16 (Mod/Dense SAV under Winter season) - fill in new H:L letter cells
$ws.Range("H16").Value = "c"
$ws.Range("I16").Value = "d"
$ws.Range("J16").Value = "c"
$ws.Range("K16").Value = "c"
$ws.Range("L16").Value = "c"

# Row 23 (Deep/Low SAV under Spring season) - update letters
$ws.Range("H23").Value = "e"
$ws.Range("I23").Value = "ef"
$ws.Range("J23").Value = "e"
$ws.Range("K23").Value = "ef"
$ws.Range("L23").Value = "f"

# Row 34 (Deep/Low SAV under Summer season) - update letters to all "g"
$ws.Range("H34").Value = "g"
$ws.Range("I34").Value = "g"
$ws.Range("J34").Value = "g"
$ws.Range("K34").Value = "g"
$ws.Range("L34").Value = "g"

# Update selection to match new active cell
$ws.Range("M34").Select()
